# Update cryptos list (prices in column D, volume % in column E)
# Generated from commit "Updated cryptos list on Mon Feb 27 07:00:44 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.441.63"
$ws.Range("E2").Value = "  +1.11%  "

$ws.Range("D3").Value = "1.638.62"
$ws.Range("E3").Value = "  +2.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.34"
$ws.Range("E6").Value = "  +1.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3761"
$ws.Range("E7").Value = "  -0.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.36"
$ws.Range("E8").Value = "  +0.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3641"
$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.262"
$ws.Range("E10").Value = "  -0.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08149"
$ws.Range("E11").Value = "  +0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.91"
$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.627"
$ws.Range("E14").Value = "  +0.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.364"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "1.640.28"
$ws.Range("E17").Value = "  +2.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.60"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06929"
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.14"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.547"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "23.452.55"
$ws.Range("E23").Value = "  +1.17%  "

$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.084"
$ws.Range("E25").Value = "  +2.64%  "

$ws.Range("E26").Value = "  +1.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.22"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.70"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.350"
$ws.Range("E29").Value = "  +2.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.31"
$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.294"
$ws.Range("E31").Value = "  -3.89%  "

$ws.Range("D32").Value = "1.823.94"
$ws.Range("E32").Value = "  +2.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.779"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9653"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02820"
$ws.Range("E35").Value = "  +4.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.34"
$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07307"
$ws.Range("E37").Value = "  -2.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2524"
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.119"
$ws.Range("E40").Value = "  +0.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.376"
$ws.Range("E41").Value = "  +1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7095"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.49"
$ws.Range("E43").Value = "  +0.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.21"
$ws.Range("E44").Value = "  +4.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6538"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("E46").Value = "  +1.15%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.019"
$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07954"
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.90"
$ws.Range("E50").Value = "  -2.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.199"
$ws.Range("E51").Value = "  -0.01%  "
